# Update cryptos list figures (price + 1h volume change) as scraped on
# Thu Apr  4 21:25:27 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.731.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.327.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("E8").Value = '  +3.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.318.84'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.178'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.580'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.31'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000275'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '639.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +11.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.860.09'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.81%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.897.58'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.52%  '
$ws.Range("E18").Value = '  +1.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.330.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.91'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.899'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.99'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.63%  '
$ws.Range("E27").Value = '  +4.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.93'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.57'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '591.27'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.941.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.96'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.40%  '
$ws.Range("E35").Value = '  -1.86%  '
$ws.Range("E36").Value = '  +2.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.76'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.26'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.128'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '32.63'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0683'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.38'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.337'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0417'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.128'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.34'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.23%  '
